$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire first column (A), shifting B:F left to A:E
$ws.Columns.Item(1).Delete()

# Fix the header text for the MODEL_CONDITION -> MODELCONDITION column (now column D after the shift)
$ws.Range("D1").Value = "MODELCONDITION"
